# Update the numeric dataset in columns B:D (rows 2-16) of Sheet1 with the
# newly generated values ("New datasets in 3 different formats").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(0.07512042055345571, 100.2386046567292, 7.258908554471414),
    @(0.1104829761145092, 61.94845412858466, 4.387980000158048),
    @(0.07781958185570878, 43.26026846297344, 1.139915756830903),
    @(0.1361248479443644, 34.5374730673197, 1.012189588113415),
    @(0.1035225351099934, 27.01537521275104, 1.348559377764242),
    @(0.1000806499181118, 24.14640749727333, 1.539107627242321),
    @(0.08747412445936389, 21.46097653845152, 1.106223817213171),
    @(0.06458179522492154, 18.38813779611064, 0.6229751598288759),
    @(0.06631390136163083, 16.40182690790442, 1.202474338343891),
    @(0.1291271735426476, 14.21161874149468, 0.7636113497023628),
    @(0.1498729653257748, 13.57640760789033, 0.351266134607137),
    @(0.1016548188108866, 12.60988953203615, 0.4325746427257349),
    @(0.1114655443483998, 11.67184835726648, 0.3594826975873466),
    @(0.08424850625946717, 11.0517883735593, 0.7862088618125245),
    @(0.1434652470131056, 9.814322806320668, 0.554749430032288)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
}
